$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Matt",     45477, "Run",      45, 4.58, 295, 3,  34, 5, 0, 0, "Sauntering Hippo", 4),
    @("Matt",     45477, "Walk",     26, 1,    20,  26, 0,  0, 0, 0, "Sauntering Hippo", 4),
    @("Eric",     45477, "Workout",  46, 0,    0,   8,  31, 7, 1, 0, "Agile Antelope",   4),
    @("Steven",   45477, "Walk",     30, 1.62, 33,  30, 0,  0, 0, 0, "Mighty Monkey",    4),
    @("Steven",   45477, "Walk",     27, 1,    243, 22, 4,  1, 0, 0, "Mighty Monkey",    4),
    @("Jeremiah", 45477, "Workout",  32, 0,    0,   15, 15, 2, 0, 0, "Agile Antelope",   4),
    @("Steven",   45477, "Walk",     36, 1.27, 46,  36, 0,  0, 0, 0, "Mighty Monkey",    4),
    @("Steven",   45478, "Workout",  34, 0,    0,   34, 0,  0, 0, 0, "Mighty Monkey",    4),
    @("Steven",   45478, "Walk",     18, 0.96, 49,  18, 0,  0, 0, 0, "Mighty Monkey",    4)
)

$startRow = 153

# Copy the date format (style) used by the existing "Date" column (B2) once,
# then paste-format-only onto the new date cells so we reuse the existing
# style record instead of Excel minting a brand-new (duplicate) style.
$ws.Cells.Item(2, 2).Copy() | Out-Null

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]

    $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

$ws.Application.ActiveWindow.ScrollRow = 137
$ws.Range("A162").Select()
